# Update "tpm" figures + target-cluster labels for the Ihh-Boc LR-pair sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Target cluster (column D) renames.
#    Row 4 was "MuSCs" -> becomes "Inflammatory-Mac"
#    Row 5 was "Neutrophils" -> becomes "MuSCs"
# ---------------------------------------------------------------------------
$ws.Range("D4").Value = "Inflammatory-Mac"
$ws.Range("D5").Value = "MuSCs"

# ---------------------------------------------------------------------------
# 2) Recomputed TPM-derived numeric values.
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("M2").Value = 2.586115
$ws.Range("N2").Value = 5.17223
$ws.Range("O2").Value = 0.1424410420966074
$ws.Range("P2").Value = 0.1066214632654476
$ws.Range("Q2").Value = 0.04020805398166667
$ws.Range("R2").Value = 0.24124832389
$ws.Range("S2").Value = 0.1424410420966074
$ws.Range("T2").Value = 0.1066214632654476

# Row 3
$ws.Range("O3").Value = 0.6698285531706168
$ws.Range("P3").Value = 0.7520806442948283
$ws.Range("S3").Value = 0.6698285531706168
$ws.Range("T3").Value = 0.7520806442948283

# Row 4
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.03546733333333333
$ws.Range("N4").Value = 0.106402
$ws.Range("O4").Value = 0.001953510930638328
$ws.Range("P4").Value = 0.002193393745902667
$ws.Range("Q4").Value = 0.0005514342762222221
$ws.Range("R4").Value = 0.004962908485999999
$ws.Range("S4").Value = 0.001953510930638328
$ws.Range("T4").Value = 0.002193393745902667

# Row 5
$ws.Range("K5").Value = 2
$ws.Range("M5").Value = 3.3707325
$ws.Range("N5").Value = 6.741465
$ws.Range("O5").Value = 0.1856571149886616
$ws.Range("P5").Value = 0.1389700115526186
$ws.Range("Q5").Value = 0.0524070253325
$ws.Range("R5").Value = 0.314442151995
$ws.Range("S5").Value = 0.1856571149886616
$ws.Range("T5").Value = 0.1389700115526186

# Row 6
$ws.Range("M6").Value = 0.002174666666666667
$ws.Range("N6").Value = 0.006524
$ws.Range("O6").Value = 0.0001197788134761043
$ws.Range("P6").Value = 0.0001344871412028815
$ws.Range("Q6").Value = 0.00003381099244444445
$ws.Range("R6").Value = 0.000304298932
$ws.Range("S6").Value = 0.0001197788134761043
$ws.Range("T6").Value = 0.0001344871412028815
